$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.69480537410311
$ws.Range("D2").Value = 4.683476307119589
$ws.Range("E2").Value = 12.87831376466378
$ws.Range("F2").Value = 23.15406831924245
$ws.Range("G2").Value = 26.86910782407855
$ws.Range("H2").Value = 13.77273033067101
$ws.Range("I2").Value = 22.05226081922623
$ws.Range("K2").Value = 13.36787194482838
$ws.Range("L2").Value = 9.381192024616887
$ws.Range("N2").Value = 16.6468168972976
$ws.Range("O2").Value = 20.71971039605553

$ws.Range("C3").Value = 10.61200474333063
$ws.Range("D3").Value = 4.625447023131954
$ws.Range("E3").Value = 12.81948528910206
$ws.Range("F3").Value = 23.16411999715816
$ws.Range("G3").Value = 26.87365475598542
$ws.Range("H3").Value = 13.81688045275357
$ws.Range("I3").Value = 22.08992778018359
$ws.Range("K3").Value = 12.86778357766004
$ws.Range("L3").Value = 9.381133550317317
$ws.Range("N3").Value = 16.65895923327601
$ws.Range("O3").Value = 20.77977924249091

$ws.Range("C4").Value = 10.5630207972595
$ws.Range("D4").Value = 4.588956125367816
$ws.Range("E4").Value = 12.78600615002856
$ws.Range("F4").Value = 23.17718098908303
$ws.Range("G4").Value = 26.88651888290319
$ws.Range("H4").Value = 13.84645801155976
$ws.Range("I4").Value = 22.11847514047845
$ws.Range("K4").Value = 12.5511165183656
$ws.Range("L4").Value = 9.382679454212511
$ws.Range("N4").Value = 16.66810331891496
$ws.Range("O4").Value = 20.82175867976434

$ws.Range("C5").Value = 10.54354495025164
$ws.Range("D5").Value = 4.573877483149601
$ws.Range("E5").Value = 12.77303842916971
$ws.Range("F5").Value = 23.18423358638692
$ws.Range("G5").Value = 26.89428774789471
$ws.Range("H5").Value = 13.85913140186689
$ws.Range("I5").Value = 22.13146846069031
$ws.Range("K5").Value = 12.41983425717264
$ws.Range("L5").Value = 9.383707920861378
$ws.Range("N5").Value = 16.67225536479254
$ws.Range("O5").Value = 20.84014355727224

$ws.Range("C6").Value = 10.54034085705231
$ws.Range("D6").Value = 4.571361346101328
$ws.Range("E6").Value = 12.77092623923469
$ws.Range("F6").Value = 23.18550907341341
$ws.Range("G6").Value = 26.89573013959771
$ws.Range("H6").Value = 13.86127325492924
$ws.Range("I6").Value = 22.13370804426938
$ws.Range("K6").Value = 12.39790538754235
$ws.Range("L6").Value = 9.383902783765057
$ws.Range("N6").Value = 16.67297055597582
$ws.Range("O6").Value = 20.84327342128399

$ws.Range("C7").Value = 10.56275614914978
$ws.Range("D7").Value = 4.588753602082551
$ws.Range("E7").Value = 12.7858285146759
$ws.Range("F7").Value = 23.17726910141665
$ws.Range("G7").Value = 26.88661343615329
$ws.Range("H7").Value = 13.8466264186039
$ws.Range("I7").Value = 22.11864487022374
$ws.Range("K7").Value = 12.54935481022172
$ws.Range("L7").Value = 9.382691710108865
$ws.Range("N7").Value = 16.6681575895155
$ws.Range("O7").Value = 20.82200145554881

$ws.Range("C8").Value = 10.66588204591354
$ws.Range("D8").Value = 4.66365125057818
$ws.Range("E8").Value = 12.85748827439169
$ws.Range("F8").Value = 23.15610329312924
$ws.Range("G8").Value = 26.86858211700283
$ws.Range("H8").Value = 13.78744045005337
$ws.Range("I8").Value = 22.06412198655935
$ws.Range("K8").Value = 13.19753750060517
$ws.Range("L8").Value = 9.380844047513753
$ws.Range("N8").Value = 16.65065378546771
$ws.Range("O8").Value = 20.73936208090371

$ws.Range("C9").Value = 10.88197369085391
$ws.Range("D9").Value = 4.803351555804436
$ws.Range("E9").Value = 13.01845030864902
$ws.Range("F9").Value = 23.16930823284588
$ws.Range("G9").Value = 26.91330689735073
$ws.Range("H9").Value = 13.69099562469708
$ws.Range("I9").Value = 22.0003119686328
$ws.Range("K9").Value = 14.38547315117995
$ws.Range("L9").Value = 9.38973357311145
$ws.Range("N9").Value = 16.62967760322063
$ws.Range("O9").Value = 20.61790695813352

$ws.Range("C10").Value = 11.04798295507567
$ws.Range("D10").Value = 4.901172831568948
$ws.Range("E10").Value = 13.14841138541506
$ws.Range("F10").Value = 23.21235428871409
$ws.Range("G10").Value = 26.99505281221298
$ws.Range("H10").Value = 13.6321333646617
$ws.Range("I10").Value = 21.97982925817091
$ws.Range("K10").Value = 15.19959284848479
$ws.Range("L10").Value = 9.403833309576491
$ws.Range("N10").Value = 16.62233693080839
$ws.Range("O10").Value = 20.55362836689721

$ws.Range("C11").Value = 11.12481162603425
$ws.Range("D11").Value = 4.944539558094471
$ws.Range("E11").Value = 13.20989813778307
$ws.Range("F11").Value = 23.23914942414825
$ws.Range("G11").Value = 27.04281742139706
$ws.Range("H11").Value = 13.6079674890567
$ws.Range("I11").Value = 21.9762580477979
$ws.Range("K11").Value = 15.55586537140003
$ws.Range("L11").Value = 9.411874992916527
$ws.Range("N11").Value = 16.62073488759155
$ws.Range("O11").Value = 20.52984318582401

$ws.Range("C12").Value = 11.15407015429273
$ws.Range("D12").Value = 4.960791404866549
$ws.Range("E12").Value = 13.23350599125292
$ws.Range("F12").Value = 23.25032876523046
$ws.Range("G12").Value = 27.0624182751392
$ws.Range("H12").Value = 13.59919253469771
$ws.Range("I12").Value = 21.97573230221438
$ws.Range("K12").Value = 15.68865899865376
$ws.Range("L12").Value = 9.415152545020614
$ws.Range("N12").Value = 16.62037668942892
$ws.Range("O12").Value = 20.52162350002449

$ws.Range("C13").Value = 11.14776183184498
$ws.Range("D13").Value = 4.957298961496993
$ws.Range("E13").Value = 13.22840746446248
$ws.Range("F13").Value = 23.24787525677081
$ws.Range("G13").Value = 27.05812971988314
$ws.Range("H13").Value = 13.60106563841203
$ws.Range("I13").Value = 21.97580876992036
$ws.Range("K13").Value = 15.66015513850093
$ws.Range("L13").Value = 9.414436361214634
$ws.Range("N13").Value = 16.62044280553114
$ws.Range("O13").Value = 20.52335870292321

$ws.Range("C14").Value = 11.12721556594714
$ws.Range("D14").Value = 4.945880059042389
$ws.Range("E14").Value = 13.21183396827727
$ws.Range("F14").Value = 23.24004849164492
$ws.Range("G14").Value = 27.04439971710978
$ws.Range("H14").Value = 13.60723802599607
$ws.Range("I14").Value = 21.97619822877068
$ws.Range("K14").Value = 15.56683335101224
$ws.Range("L14").Value = 9.412139993064113
$ws.Range("N14").Value = 16.62070044681849
$ws.Range("O14").Value = 20.52915115507748

$ws.Range("C15").Value = 11.1146512037584
$ws.Range("D15").Value = 4.938863275592823
$ws.Range("E15").Value = 13.20172394771679
$ws.Range("F15").Value = 23.23538868714089
$ws.Range("G15").Value = 27.03618651579401
$ws.Range("H15").Value = 13.61106779623555
$ws.Range("I15").Value = 21.97654442842228
$ws.Range("K15").Value = 15.50939249652142
$ws.Range("L15").Value = 9.410763605152663
$ws.Range("N15").Value = 16.62089057545234
$ws.Range("O15").Value = 20.5328017983745

$ws.Range("C16").Value = 11.04298640394406
$ws.Range("D16").Value = 4.898315289784003
$ws.Range("E16").Value = 13.14443933702455
$ws.Range("F16").Value = 23.21074794507751
$ws.Range("G16").Value = 26.99214354249612
$ws.Range("H16").Value = 13.63376525906487
$ws.Range("I16").Value = 21.98017831405907
$ws.Range("K16").Value = 15.17601795300266
$ws.Range("L16").Value = 9.403340376651247
$ws.Range("N16").Value = 16.62247647147889
$ws.Range("O16").Value = 20.5552928051776

$ws.Range("C17").Value = 10.99934148138158
$ws.Range("D17").Value = 4.873145012835741
$ws.Range("E17").Value = 13.10989194724494
$ws.Range("F17").Value = 23.19747624407636
$ws.Range("G17").Value = 26.9678292491562
$ws.Range("H17").Value = 13.64835862751999
$ws.Range("I17").Value = 21.98387976023107
$ws.Range("K17").Value = 14.96782673608579
$ws.Range("L17").Value = 9.399202182367535
$ws.Range("N17").Value = 16.62389337267457
$ws.Range("O17").Value = 20.57048964625324

$ws.Range("C18").Value = 10.97436280610071
$ws.Range("D18").Value = 4.858561542748249
$ws.Range("E18").Value = 13.09024498546285
$ws.Range("F18").Value = 23.19052202163604
$ws.Range("G18").Value = 26.95484058510819
$ws.Range("H18").Value = 13.65699805345164
$ws.Range("I18").Value = 21.98654964026664
$ws.Range("K18").Value = 14.84676069221582
$ws.Range("L18").Value = 9.396975317466492
$ws.Range("N18").Value = 16.62487194039789
$ws.Range("O18").Value = 20.57974381673455

$ws.Range("C19").Value = 10.96592758316524
$ws.Range("D19").Value = 4.853605819454268
$ws.Range("E19").Value = 13.08363178048476
$ws.Range("F19").Value = 23.1882842337243
$ws.Range("G19").Value = 26.95061414885019
$ws.Range("H19").Value = 13.65996539321113
$ws.Range("I19").Value = 21.98754649973832
$ws.Range("K19").Value = 14.80554621945407
$ws.Range("L19").Value = 9.396247720561282
$ws.Range("N19").Value = 16.62523140625781
$ws.Range("O19").Value = 20.58296519906146

$ws.Range("C20").Value = 11.00397481657546
$ws.Range("D20").Value = 4.875835476241575
$ws.Range("E20").Value = 13.11354653061031
$ws.Range("F20").Value = 23.19881876291752
$ws.Range("G20").Value = 26.97031448381765
$ws.Range("H20").Value = 13.64677970368895
$ws.Range("I20").Value = 21.98342974952392
$ws.Range("K20").Value = 14.99012638807738
$ws.Range("L20").Value = 9.399626844532774
$ws.Range("N20").Value = 16.62372561843223
$ws.Range("O20").Value = 20.56881876586829

$ws.Range("C21").Value = 11.1332462052903
$ws.Range("D21").Value = 4.949238743260449
$ws.Range("E21").Value = 13.21669334019151
$ws.Range("F21").Value = 23.24231942127485
$ws.Range("G21").Value = 27.04839155156031
$ws.Range("H21").Value = 13.60541483241204
$ws.Range("I21").Value = 21.9760614025185
$ws.Range("K21").Value = 15.59430242676369
$ws.Range("L21").Value = 9.412808200891854
$ws.Range("N21").Value = 16.62061803920606
$ws.Range("O21").Value = 20.52742838543492

$ws.Range("C22").Value = 11.21868445801001
$ws.Range("D22").Value = 4.996216834618183
$ws.Range("E22").Value = 13.28598727567465
$ws.Range("F22").Value = 23.27676534556087
$ws.Range("G22").Value = 27.10823549594273
$ws.Range("H22").Value = 13.58057322518616
$ws.Range("I22").Value = 21.97606349922362
$ws.Range("K22").Value = 15.97677966561909
$ws.Range("L22").Value = 9.422776324176517
$ws.Range("N22").Value = 16.62003461462979
$ws.Range("O22").Value = 20.50496710531166

$ws.Range("C23").Value = 11.17300506094033
$ws.Range("D23").Value = 4.971237142251252
$ws.Range("E23").Value = 13.24883712758759
$ws.Range("F23").Value = 23.25783234354445
$ws.Range("G23").Value = 27.07549218662508
$ws.Range("H23").Value = 13.59363079952145
$ws.Range("I23").Value = 21.97562162011336
$ws.Range("K23").Value = 15.77380591227956
$ws.Range("L23").Value = 9.417332928644036
$ws.Range("N23").Value = 16.62021400756669
$ws.Range("O23").Value = 20.51653432448374

$ws.Range("C24").Value = 11.00187973083761
$ws.Range("D24").Value = 4.874619468344122
$ws.Range("E24").Value = 13.11189362320312
$ws.Range("F24").Value = 23.19820970456718
$ws.Range("G24").Value = 26.96918782505697
$ws.Range("H24").Value = 13.64749275800586
$ws.Range("I24").Value = 21.98363151149853
$ws.Range("K24").Value = 14.9800489881849
$ws.Range("L24").Value = 9.399434380436908
$ws.Range("N24").Value = 16.62380094925661
$ws.Range("O24").Value = 20.56957255972128

$ws.Range("C25").Value = 10.82215319607234
$ws.Range("D25").Value = 4.766373644501865
$ws.Range("E25").Value = 12.97279286616389
$ws.Range("F25").Value = 23.15987468816792
$ws.Range("G25").Value = 26.89261740604178
$ws.Range("H25").Value = 13.71498242956541
$ws.Range("I25").Value = 22.01294491292411
$ws.Range("K25").Value = 14.07389518830833
$ws.Range("L25").Value = 9.385994606297363
$ws.Range("N25").Value = 16.63392991595476
$ws.Range("O25").Value = 20.64639518891132
